$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: 3/11/15, 2:28 PM - 5:00 PM -------------------------------
# Copy time-format styles already used on A/B from row 10 isn't needed
# (A11/B11 already carry style s="7"); just copy the date style from an
# existing date cell so we reuse the existing numFmt (style s="4") rather
# than create a brand-new number format entry.
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("A11").Value2 = 0.60277777777777775
$ws.Range("B11").Value2 = 0.70833333333333337
$ws.Range("C11").Value2 = 42074
$ws.Range("E11").Value2 = "Updated the website once again to make sure that the website is splite into boxes. The layout of these boxes are set to change according to the screen size. This is similar to how a website readjusts itself."
$ws.Rows.Item(11).RowHeight = 60

# --- Row 12: 3/16/15, 3:00 PM - 5:00 PM -------------------------------
$ws.Range("C10").Copy()
$ws.Range("C12").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("E12").PasteSpecial(-4122)

$ws.Range("A12").Value2 = 0.625
$ws.Range("B12").Value2 = 0.70833333333333337
$ws.Range("C12").Value2 = 42079
$ws.Range("E12").Value2 = "Changed the font color of a title, subtitle, etc. on the website to red when you hover over it. We used the hover, mouseover, and mouseout functions to accomplish this task."
$ws.Rows.Item(12).RowHeight = 45

# --- Row 13: 3/23/15, 2:30 PM - 5:00 PM -------------------------------
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("A13").Value2 = 0.60416666666666663
$ws.Range("B13").Value2 = 0.70833333333333337
$ws.Range("C13").Value2 = 42086
$ws.Range("E13").Value2 = "Tried to implement a way for text to shake using Jquery ui. "

# --- Move the view / selection to the newly added entries -------------
$ws.Range("E13").Select()
